{"js": "// Applies the LOB1202.docx content update:\n//  - bumps the \"Ativa\u00e7\u00e3o\" date from 01/01/2018 to 01/01/2024\n//  - appends a sentence about didactic trips to both the PT and EN\n//    \"Programa\" paragraphs\n//  - rewrites the \"M\u00e9todo\", \"Crit\u00e9rio\" and \"Norma de recupera\u00e7\u00e3o\"\n//    fields under \"Avalia\u00e7\u00e3o\"\n\nasync function replaceText(body, oldText, newText) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1) Activation date\nawait replaceText(body, \"Ativa\u00e7\u00e3o: 01/01/2018\", \"Ativa\u00e7\u00e3o: 01/01/2024\");\n\n// 2) Portuguese \"Programa\" paragraph - append sentence about field trips\nawait replaceText(\n  body,\n  \"Defini\u00e7\u00e3o de projeto e seus principais atributos e caracter\u00edsticas; conceitos do PMBoK (Project managment body of knowledge). Planejamento estrat\u00e9gico. Desenvolvimento Sustent\u00e1vel: O que \u00e9 desenvolvimento sustent\u00e1vel? Conv\u00eanios, tratados e pol\u00edticas de alcance internacional realizado em torno do desenvolvimento sustent\u00e1vel. Os desafios do desenvolvimento sustent\u00e1vel. Processos e metodologia do gerenciamento de projetos ambientais. Ferramentas de planejamento, monitoramento e controle. Estudo dos riscos e problemas comuns na gest\u00e3o de projetos ambientais. An\u00e1lise de casos reais envolvendo sele\u00e7\u00e3o, administra\u00e7\u00e3o e desenvolvimento de projetos aplicados \u00e0 gest\u00e3o ambiental\",\n  \"Defini\u00e7\u00e3o de projeto e seus principais atributos e caracter\u00edsticas; conceitos do PMBoK (Project managment body of knowledge). Planejamento estrat\u00e9gico. Desenvolvimento Sustent\u00e1vel: O que \u00e9 desenvolvimento sustent\u00e1vel? Conv\u00eanios, tratados e pol\u00edticas de alcance internacional realizado em torno do desenvolvimento sustent\u00e1vel. Os desafios do desenvolvimento sustent\u00e1vel. Processos e metodologia do gerenciamento de projetos ambientais. Ferramentas de planejamento, monitoramento e controle. Estudo dos riscos e problemas comuns na gest\u00e3o de projetos ambientais. An\u00e1lise de casos reais envolvendo sele\u00e7\u00e3o, administra\u00e7\u00e3o e desenvolvimento de projetos aplicados \u00e0 gest\u00e3o ambiental. A disciplina pode contar com viagens did\u00e1ticas para complementa\u00e7\u00e3o do conte\u00fado da disciplina.\"\n);\n\n// 3) English \"Programa\" paragraph - append sentence about field trips\nawait replaceText(\n  body,\n  \"Definition of project and its main attributes and characteristics. Project Management Body of Knowledge concepts. Strategic planning. Sustainable Development: Definition. Agreements, treaties and policies international carried out around sustainable development. The challenges of sustainable. Processes and Methodology of management of environmental projects. Planning tools, monitoring and control. Study of the risks and common problems in the management of environmental projects. Analysis of cases involving selection, management and development of environmental management applied projects.\",\n  \"Definition of project and its main attributes and characteristics. Project Management Body of Knowledge concepts. Strategic planning. Sustainable Development: Definition. Agreements, treaties and policies international carried out around sustainable development. The challenges of sustainable. Processes and Methodology of management of environmental projects. Planning tools, monitoring and control. Study of the risks and common problems in the management of environmental projects. Analysis of cases involving selection, management and development of environmental management applied projects. The discipline may have didactic trips to complement the content of the discipline.\"\n);\n\n// 4) \"M\u00e9todo:\" evaluation text\nawait replaceText(\n  body,\n  \"Avalia\u00e7\u00e3o composta por 2 (duas) provas e um projetoPara os alunos que perderem uma das provas ser\u00e1 oferecida uma substitutiva no final do semestre, que incluir\u00e1 toda a mat\u00e9ria da disciplina.\",\n  \"Aulas te\u00f3ricas e pr\u00e1ticas, exerc\u00edcios dirigidos. Avalia\u00e7\u00e3o baseada em provas, exerc\u00edcios e trabalhos pr\u00e1ticos e relat\u00f3rios\"\n);\n\n// 5) \"Crit\u00e9rio:\" evaluation text\nawait replaceText(\n  body,\n  \"Nota final = (nota prova 1 + nota da prova 2 + nota do projeto)/3.\",\n  \"M\u00e9dia das avalia\u00e7\u00f5es aplicadas\"\n);\n\n// 6) \"Norma de recupera\u00e7\u00e3o:\" evaluation text\nawait replaceText(\n  body,\n  \"Prova \u00fanica com todo o conte\u00fado da disciplina, sendo que a nota [(nota final do semestre + nota recupera\u00e7\u00e3o)/2] dever\u00e1 ser igual ou superior a 5,0 (cinco).\",\n  \"A Nota Final ser\u00e1 composta pela M\u00e9dia obtida da Nota do Per\u00edodo somada \u00e0 Nota de Recupera\u00e7\u00e3o e dividido por dois\"\n);\n", "ps1": "# Applies the LOB1202.docx content update:\n#  - bumps the \"Ativa\u00e7\u00e3o\" date from 01/01/2018 to 01/01/2024\n#  - appends a sentence about didactic trips to both the PT and EN\n#    \"Programa\" paragraphs\n#  - rewrites the \"M\u00e9todo\", \"Crit\u00e9rio\" and \"Norma de recupera\u00e7\u00e3o\"\n#    fields under \"Avalia\u00e7\u00e3o\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute(\n        $oldText,    # FindText\n        $false,      # MatchCase\n        $true,       # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $newText,    # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n    if (-not $result) {\n        throw \"Find/Replace failed for: $oldText\"\n    }\n}\n\n# 1) Activation date\nReplace-Text \"Ativa\u00e7\u00e3o: 01/01/2018\" \"Ativa\u00e7\u00e3o: 01/01/2024\"\n\n# 2) Portuguese \"Programa\" paragraph - append sentence about field trips\nReplace-Text `\n    \"Defini\u00e7\u00e3o de projeto e seus principais atributos e caracter\u00edsticas; conceitos do PMBoK (Project managment body of knowledge). Planejamento estrat\u00e9gico. Desenvolvimento Sustent\u00e1vel: O que \u00e9 desenvolvimento sustent\u00e1vel? Conv\u00eanios, tratados e pol\u00edticas de alcance internacional realizado em torno do desenvolvimento sustent\u00e1vel. Os desafios do desenvolvimento sustent\u00e1vel. Processos e metodologia do gerenciamento de projetos ambientais. Ferramentas de planejamento, monitoramento e controle. Estudo dos riscos e problemas comuns na gest\u00e3o de projetos ambientais. An\u00e1lise de casos reais envolvendo sele\u00e7\u00e3o, administra\u00e7\u00e3o e desenvolvimento de projetos aplicados \u00e0 gest\u00e3o ambiental\" `\n    \"Defini\u00e7\u00e3o de projeto e seus principais atributos e caracter\u00edsticas; conceitos do PMBoK (Project managment body of knowledge). Planejamento estrat\u00e9gico. Desenvolvimento Sustent\u00e1vel: O que \u00e9 desenvolvimento sustent\u00e1vel? Conv\u00eanios, tratados e pol\u00edticas de alcance internacional realizado em torno do desenvolvimento sustent\u00e1vel. Os desafios do desenvolvimento sustent\u00e1vel. Processos e metodologia do gerenciamento de projetos ambientais. Ferramentas de planejamento, monitoramento e controle. Estudo dos riscos e problemas comuns na gest\u00e3o de projetos ambientais. An\u00e1lise de casos reais envolvendo sele\u00e7\u00e3o, administra\u00e7\u00e3o e desenvolvimento de projetos aplicados \u00e0 gest\u00e3o ambiental. A disciplina pode contar com viagens did\u00e1ticas para complementa\u00e7\u00e3o do conte\u00fado da disciplina.\"\n\n# 3) English \"Programa\" paragraph - append sentence about field trips\nReplace-Text `\n    \"Definition of project and its main attributes and characteristics. Project Management Body of Knowledge concepts. Strategic planning. Sustainable Development: Definition. Agreements, treaties and policies international carried out around sustainable development. The challenges of sustainable. Processes and Methodology of management of environmental projects. Planning tools, monitoring and control. Study of the risks and common problems in the management of environmental projects. Analysis of cases involving selection, management and development of environmental management applied projects.\" `\n    \"Definition of project and its main attributes and characteristics. Project Management Body of Knowledge concepts. Strategic planning. Sustainable Development: Definition. Agreements, treaties and policies international carried out around sustainable development. The challenges of sustainable. Processes and Methodology of management of environmental projects. Planning tools, monitoring and control. Study of the risks and common problems in the management of environmental projects. Analysis of cases involving selection, management and development of environmental management applied projects. The discipline may have didactic trips to complement the content of the discipline.\"\n\n# 4) \"M\u00e9todo:\" evaluation text\nReplace-Text `\n    \"Avalia\u00e7\u00e3o composta por 2 (duas) provas e um projetoPara os alunos que perderem uma das provas ser\u00e1 oferecida uma substitutiva no final do semestre, que incluir\u00e1 toda a mat\u00e9ria da disciplina.\" `\n    \"Aulas te\u00f3ricas e pr\u00e1ticas, exerc\u00edcios dirigidos. Avalia\u00e7\u00e3o baseada em provas, exerc\u00edcios e trabalhos pr\u00e1ticos e relat\u00f3rios\"\n\n# 5) \"Crit\u00e9rio:\" evaluation text\nReplace-Text `\n    \"Nota final = (nota prova 1 + nota da prova 2 + nota do projeto)/3.\" `\n    \"M\u00e9dia das avalia\u00e7\u00f5es aplicadas\"\n\n# 6) \"Norma de recupera\u00e7\u00e3o:\" evaluation text\nReplace-Text `\n    \"Prova \u00fanica com todo o conte\u00fado da disciplina, sendo que a nota [(nota final do semestre + nota recupera\u00e7\u00e3o)/2] dever\u00e1 ser igual ou superior a 5,0 (cinco).\" `\n    \"A Nota Final ser\u00e1 composta pela M\u00e9dia obtida da Nota do Per\u00edodo somada \u00e0 Nota de Recupera\u00e7\u00e3o e dividido por dois\"\n\nWrite-Output \"LOB1202 updates applied\"\n"}
